$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AB1:AC1").EntireColumn.Insert()
$ws.Range("AC1").Value = "利害關係人"
$ws.Range("AB1").Value = "特定資產記號"
